$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D2").Value = 87
$ws.Range("H2").Value = 99
$ws.Range("J2").Value = 111
$ws.Range("G3").Value = 133
$ws.Range("J3").Value = 211
$ws.Range("J6").Value = 19
$ws.Range("B9").Value = 354
$ws.Range("C9").Value = 454
$ws.Range("E9").Value = 428
$ws.Range("F9").Value = 483
$ws.Range("H9").Value = 421
$ws.Range("I9").Value = 475
$ws.Range("K9").Value = 466
$ws.Range("B10").Value = 1243
$ws.Range("D10").Value = 1680
$ws.Range("E10").Value = 1991
$ws.Range("F10").Value = 1975
$ws.Range("G10").Value = 857
$ws.Range("H10").Value = 558
$ws.Range("I10").Value = 798
$ws.Range("J10").Value = 677
$ws.Range("K10").Value = 625
$ws.Range("B11").Value = 1722
$ws.Range("C11").Value = 2063
$ws.Range("D11").Value = 2286
$ws.Range("E11").Value = 2634
$ws.Range("F11").Value = 2672
$ws.Range("G11").Value = 1498
$ws.Range("H11").Value = 1238
$ws.Range("I11").Value = 1595
$ws.Range("J11").Value = 1411
$ws.Range("K11").Value = 1458

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 13
$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 22
$ws.Range("E8").Value = 102
$ws.Range("I8").Value = 80
$ws.Range("F19").Value = 50
$ws.Range("E21").Value = 24
$ws.Range("H21").Value = 20
$ws.Range("B27").Value = 20
$ws.Range("B28").Value = 93
$ws.Range("C28").Value = 124
$ws.Range("F28").Value = 115
$ws.Range("J28").Value = 59
$ws.Range("K28").Value = 85
$ws.Range("G29").Value = 13
$ws.Range("E32").Value = 137
$ws.Range("C36").Value = 93
$ws.Range("D36").Value = 75
$ws.Range("F38").Value = 11
$ws.Range("I47").Value = 48
$ws.Range("J47").Value = 38
$ws.Range("B53").Value = 242
$ws.Range("E53").Value = 659
$ws.Range("G53").Value = 237
$ws.Range("H53").Value = 185
$ws.Range("I53").Value = 301
$ws.Range("D54").Value = 17
$ws.Range("B63").Value = 13
$ws.Range("B65").Value = 30
$ws.Range("F67").Value = 19
$ws.Range("J68").Value = 13
$ws.Range("H70").Value = 25
$ws.Range("I72").Value = 15
$ws.Range("E74").Value = 73
$ws.Range("E77").Value = 69
$ws.Range("G77").Value = 51
$ws.Range("E78").Value = 46
$ws.Range("J80").Value = 10
$ws.Range("H89").Value = 11
$ws.Range("H92").Value = 23
$ws.Range("E95").Value = 95
$ws.Range("H95").Value = 18
$ws.Range("B97").Value = 29
$ws.Range("K97").Value = 13
$ws.Range("B99").Value = 1722
$ws.Range("C99").Value = 2063
$ws.Range("D99").Value = 2286
$ws.Range("E99").Value = 2634
$ws.Range("F99").Value = 2672
$ws.Range("G99").Value = 1498
$ws.Range("H99").Value = 1238
$ws.Range("I99").Value = 1595
$ws.Range("J99").Value = 1411
$ws.Range("K99").Value = 1458

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("B6").Value = 17
$ws.Range("B7").Value = 20

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("E9").Value = 46
$ws.Range("G9").Value = 28
$ws.Range("E10").Value = 69
$ws.Range("G10").Value = 51

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("D6").Value = 13
$ws.Range("D7").Value = 22

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E8").Value = 56
$ws.Range("I8").Value = 39
$ws.Range("E9").Value = 102
$ws.Range("I9").Value = 80

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("E8").Value = 15
$ws.Range("H8").Value = 8
$ws.Range("E9").Value = 24
$ws.Range("H9").Value = 20

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E8").Value = 78
$ws.Range("E9").Value = 137

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("F9").Value = 29
$ws.Range("F10").Value = 50

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("D2").Value = 7
$ws.Range("C8").Value = 31
$ws.Range("C10").Value = 93
$ws.Range("D10").Value = 75

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("B8").Value = 29
$ws.Range("B9").Value = 195
$ws.Range("E9").Value = 582
$ws.Range("G9").Value = 157
$ws.Range("H9").Value = 97
$ws.Range("I9").Value = 182
$ws.Range("B10").Value = 242
$ws.Range("E10").Value = 659
$ws.Range("G10").Value = 237
$ws.Range("H10").Value = 185
$ws.Range("I10").Value = 301

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("H7").Value = 12
$ws.Range("H8").Value = 25

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("B6").Value = 12
$ws.Range("B8").Value = 30

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("H4").Value = 4
$ws.Range("H6").Value = 11

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("E5").Value = 39
$ws.Range("E6").Value = 46

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 10
$ws.Range("C7").Value = 41
$ws.Range("K7").Value = 38
$ws.Range("B8").Value = 60
$ws.Range("F8").Value = 64
$ws.Range("B9").Value = 93
$ws.Range("C9").Value = 124
$ws.Range("F9").Value = 115
$ws.Range("J9").Value = 59
$ws.Range("K9").Value = 85

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I7").Value = 26
$ws.Range("J7").Value = 24
$ws.Range("I8").Value = 48
$ws.Range("J8").Value = 38

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("G3").Value = 2
$ws.Range("G9").Value = 13

$ws = $wb.Worksheets.Item("River North")
$ws.Range("E6").Value = 67
$ws.Range("E7").Value = 73

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("H2").Value = 4
$ws.Range("H9").Value = 23

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("F7").Value = 17
$ws.Range("F8").Value = 19

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("B6").Value = 14
$ws.Range("K6").Value = 4
$ws.Range("B7").Value = 29
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J4").Value = 1
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("D5").Value = 11
$ws.Range("D6").Value = 17

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("E5").Value = 8
$ws.Range("H6").Value = 11
$ws.Range("E7").Value = 95
$ws.Range("H7").Value = 18

$ws = $wb.Worksheets.Item("New City")
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 13

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("I4").Value = 10
$ws.Range("I6").Value = 15

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("D5").Value = 7
$ws.Range("D6").Value = 8

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 11
